$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Insert the new sheet "Capacità_trans_ MW_monodir" right after "Distanza km"
# (i.e. right before "Raw data - tr. interna 2023"), matching the diff that
# moves all "Raw data - ..." sheets one slot to the right.
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("Raw data - tr. interna 2023")
$ws = $wb.Worksheets.Add($refSheet)
$ws.Name = "Capacità_trans_ MW_monodir"

# Grab the exact cell style (bold "Aptos Narrow", thin box border, centered)
# already used for header/label cells elsewhere in the workbook, and stamp it
# onto the new sheet's header row + label column before writing values.
$styleSrc = $wb.Worksheets.Item("Capacità di trasmissione MW").Range("A2")
$styleSrc.Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$ws.Range("A2:A9").PasteSpecial(-4122)

# Header row (zone names across the top)
$ws.Range("B1").Value = "NORD"
$ws.Range("C1").Value = "CNOR"
$ws.Range("D1").Value = "CSUD"
$ws.Range("E1").Value = "SUD"
$ws.Range("F1").Value = "CALA"
$ws.Range("G1").Value = "SICI"
$ws.Range("H1").Value = "SARD"

# Row labels (zone names down the left side)
$ws.Range("A2").Value = "NORD"
$ws.Range("A3").Value = "CNOR"
$ws.Range("A4").Value = "CSUD"
$ws.Range("A5").Value = "SUD"
$ws.Range("A6").Value = "CALA"
$ws.Range("A7").Value = "SICI"
$ws.Range("A8").Value = "SARD"

# Symmetrised ("monodirezionale") transmission-capacity matrix [MW]
$matrix = @{
    "B2" = 0;    "C2" = 4300; "D2" = 0;    "E2" = 0;    "F2" = 0;    "G2" = 0;    "H2" = 0;
    "B3" = 4300; "C3" = 0;    "D3" = 2900; "E3" = 0;    "F3" = 0;    "G3" = 0;    "H3" = 300;
    "B4" = 0;    "C4" = 2900; "D4" = 0;    "E4" = 5200; "F4" = 0;    "G4" = 0;    "H4" = 900;
    "B5" = 0;    "C5" = 0;    "D5" = 5200; "E5" = 0;    "F5" = 2400; "G5" = 0;    "H5" = 0;
    "B6" = 0;    "C6" = 0;    "D6" = 0;    "E6" = 2400; "F6" = 0;    "G6" = 1600; "H6" = 1100;
    "B7" = 0;    "C7" = 0;    "D7" = 0;    "E7" = 0;    "F7" = 1600; "G7" = 0;    "H7" = 0;
    "B8" = 0;    "C8" = 300;  "D8" = 900;  "E8" = 0;    "F8" = 1100; "G8" = 0;    "H8" = 0;
}
foreach ($addr in $matrix.Keys) {
    $ws.Range($addr).Value = $matrix[$addr]
}

# Row 9 only carries the (empty) label-column style, matching the source diff.

# Make this the active sheet/tab with F9 selected, like the captured diff.
$ws.Activate()
$ws.Range("F9").Select()
